$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix the "JCL Insight 2024" typo -> "JCI Insight 2024" (row 4, Paper column) ---
$ws.Range("C4").Value = "JCI Insight 2024"

# --- Row 4 (Sparsentan / FSGS paper) now points to a new pubmed URL instead of the
#     insight.jci.org link, and that cell no longer carries a hyperlink ---
$ws.Range("D4").Value = "https://pubmed.ncbi.nlm.nih.gov/39226116/"

# --- Hyperlinks: only D2 and D3 keep live hyperlinks; D4's hyperlink is removed.
#     Remove every hyperlink on the sheet and re-add just the two that should remain,
#     since this host's Hyperlinks collection is sheet-scoped. ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D3"), "https://pubmed.ncbi.nlm.nih.gov/38996810/")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://pubmed.ncbi.nlm.nih.gov/35522041/")

# --- Update the active selection left by the author: B10 -> D4 ---
$ws.Range("D4").Select()
